$d = $word.ActiveDocument

# --- Change 1 ------------------------------------------------------------
# The very first paragraph in the body is empty (just a paragraph mark) but
# carries run-level formatting (Times New Roman, bold, 28pt) on that mark.
# The target adds a <w:vertAlign w:val="subscript"/> to that mark's rPr,
# i.e. turns the (invisible) paragraph mark into subscript.
#
# A Range that spans *only* the paragraph mark of an empty paragraph won't
# pick up Font property writes directly in this host, so we briefly type a
# placeholder character, format it (which also marks the paragraph mark as
# subscript), then delete the placeholder again - the mark keeps the
# formatting that was applied while the placeholder existed.
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertBefore("X")
$d.Paragraphs.Item(1).Range.Font.Subscript = $true
$d.Range(0, 1).Delete()

# --- Change 2 ------------------------------------------------------------
# "IMPORTANT: You must give this evaluation grid to your teacher with page
# 1." used to be typed as several runs ("IMPORTANT", ": ", "Y",
# "ou must give..."). Collapse the three runs that follow "IMPORTANT" (": ",
# "Y", "ou must give this evaluation grid to your teacher with page 1.")
# into a single run with the same formatting.
$d.Content.Find.Execute( `
    ": You must give this evaluation grid to your teacher with page 1.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    ": You must give this evaluation grid to your teacher with page 1.", 2)
